$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing bordered/bold header style (from A1:B1) across C1:D1
# by copying formats only, so C1/D1 pick up the same cellXf as A1/B1.
$ws.Range("A1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)

# Header row (row 1)
$ws.Range("A1").Value = "MIGRATION DATE"
$ws.Range("B1").Value = "FINANCIAL INSTITUTION NAME"
$ws.Range("C1").Value = "ENTITY ID"
$ws.Range("D1").Value = "ADDRESS"

# Data row (row 2) - plain unstyled cells.
# "2025-10-17" must stay literal text, not get auto-converted to a date
# serial, so force text format before assigning, then clear the format
# back off so no stray number-format style is left on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-10-17"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = "ZZZ"
$ws.Range("C2").Value = "456CDX009"
$ws.Range("D2").Value = "Anna Nagar"
